# Rebuild the workbook to reflect the corrected forecast output:
# Sheet1 gains an "Order Week" column and its dates/values are recomputed,
# and three new summary/forecast sheets are added.

$wb = $excel.ActiveWorkbook

# --- Sheet1: "Sales vs PO" ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sales vs PO"

# Insert a new column C ("Order Week"), shifting the old PO_Requested_Qty to D
$ws1.Columns.Item(3).Insert()
$ws1.Cells.Item(1,3).Value = "Order Week"

# Give column C the same date style as column A (copy number format from A2)
$ws1.Range("A2").Copy($ws1.Range("C2:C14"))

# Update column A (ds, now +6 days), C (Order Week, former ds values) and D (reset to 0)
$ws1.Cells.Item(2,1).Value = 45571
$ws1.Cells.Item(2,3).Value = 45565
$ws1.Cells.Item(2,4).Value = 0
$ws1.Cells.Item(3,1).Value = 45578
$ws1.Cells.Item(3,3).Value = 45572
$ws1.Cells.Item(3,4).Value = 0
$ws1.Cells.Item(4,1).Value = 45585
$ws1.Cells.Item(4,3).Value = 45579
$ws1.Cells.Item(4,4).Value = 0
$ws1.Cells.Item(5,1).Value = 45592
$ws1.Cells.Item(5,3).Value = 45586
$ws1.Cells.Item(5,4).Value = 0
$ws1.Cells.Item(6,1).Value = 45599
$ws1.Cells.Item(6,3).Value = 45593
$ws1.Cells.Item(6,4).Value = 0
$ws1.Cells.Item(7,1).Value = 45606
$ws1.Cells.Item(7,3).Value = 45600
$ws1.Cells.Item(7,4).Value = 0
$ws1.Cells.Item(8,1).Value = 45613
$ws1.Cells.Item(8,3).Value = 45607
$ws1.Cells.Item(8,4).Value = 0
$ws1.Cells.Item(9,1).Value = 45620
$ws1.Cells.Item(9,3).Value = 45614
$ws1.Cells.Item(9,4).Value = 0
$ws1.Cells.Item(10,1).Value = 45627
$ws1.Cells.Item(10,3).Value = 45621
$ws1.Cells.Item(10,4).Value = 0
$ws1.Cells.Item(11,1).Value = 45634
$ws1.Cells.Item(11,3).Value = 45628
$ws1.Cells.Item(11,4).Value = 0
$ws1.Cells.Item(12,1).Value = 45641
$ws1.Cells.Item(12,3).Value = 45635
$ws1.Cells.Item(12,4).Value = 0
$ws1.Cells.Item(13,1).Value = 45648
$ws1.Cells.Item(13,3).Value = 45642
$ws1.Cells.Item(13,4).Value = 0
$ws1.Cells.Item(14,1).Value = 45655
$ws1.Cells.Item(14,3).Value = 45649
$ws1.Cells.Item(14,4).Value = 0

# --- Sheet2: "Weekly Growth" ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Weekly Growth"

$ws1.Range("A1").Copy($ws2.Range("A1"))
$ws2.Range("A1").Value = "ds"
$ws1.Range("A1").Copy($ws2.Range("B1"))
$ws2.Range("B1").Value = "PO_Requested_Qty"
$ws1.Range("A1").Copy($ws2.Range("C1"))
$ws2.Range("C1").Value = "Growth%"

$ws1.Range("A2").Copy($ws2.Range("A2:A3"))
$ws2.Cells.Item(2,1).Value = 45572
$ws2.Cells.Item(2,2).Value = 120
$ws2.Cells.Item(2,3).Value = 0
$ws2.Cells.Item(3,1).Value = 45586
$ws2.Cells.Item(3,2).Value = 24
$ws2.Cells.Item(3,3).Value = -80

# --- Sheet3: "Volume Insights" ---
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Volume Insights"

$ws1.Range("A1").Copy($ws3.Range("A1"))
$ws3.Range("A1").Value = "Total_PO_Quantity"
$ws1.Range("A1").Copy($ws3.Range("B1"))
$ws3.Range("B1").Value = "Average_PO_Quantity"
$ws1.Range("A1").Copy($ws3.Range("C1"))
$ws3.Range("C1").Value = "Max_PO_Quantity"
$ws1.Range("A1").Copy($ws3.Range("D1"))
$ws3.Range("D1").Value = "Min_PO_Quantity"

$ws3.Cells.Item(2,1).Value = 144
$ws3.Cells.Item(2,2).Value = 72
$ws3.Cells.Item(2,3).Value = 120
$ws3.Cells.Item(2,4).Value = 24

# --- Sheet4: "Prediction Info" ---
$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "Prediction Info"

$ws1.Range("A1").Copy($ws4.Range("A1"))
$ws4.Range("A1").Value = "Predicted_Next_Week_PO_Quantity"
$ws4.Cells.Item(2,1).Value = 0

# Leave the first sheet active/selected
$ws1.Select()
